# Re-pull / push corrected "dSF" (column F) values after re-pulling data.
# In the prior pull, column F (dSF) was a straight copy of column E (dS0).
# This update pushes the real dSF values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 2
    7  = 3
    13 = 1
    18 = 0
    23 = 2
    26 = 1
    27 = -2
    29 = -4
    37 = 2
    39 = -2
    44 = -3
    46 = -3
    48 = -5
    50 = -5
    52 = -1
    53 = -3
    54 = -3
    61 = -2
    63 = -6
    71 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
